$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header for the new "Cap. Percent (%)" column
$ws.Range("G1").Value = "Cap. Percent (%)"

# Sum of D2:D9 (total abundance*xn) in D10, next to the existing C10 sum
$ws.Range("D10").Formula = "=SUM(D2:D9)"

# Percent-of-capture formulas in G2:G9 (absolute ref to the D10 total)
$ws.Range("G2").Formula = "=(D2/`$D`$10)*100"
$ws.Range("G3:G9").Formula = "=(D3/`$D`$10)*100"

# Match column G's width to column F's
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(6).ColumnWidth

# Move the active selection to G16, matching the recorded sheet view
$ws.Range("G16").Select()
